$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values updated for rows 2 and 3 (same new values in both rows)
$updates = @{
    "D"  = -0.0307
    "K"  = -159.1
    "L"  = -0.3867282450170151
    "M"  = 128.6
    "N"  = 0.07823336172283733
    "O"  = -0.8082966687617851
    "P"  = 128.6
    "Q"  = 0.07823336172283733
    "R"  = -0.8082966687617851
    "U"  = 1917.1
    "V"  = 1.166261102323884
    "W"  = -0.1297822008320418
    "X"  = 0.06921934836785093
    "Y"  = -0.1990015491998927
    "Z"  = -0.5273682861171645
    "AB" = 0.07017433399895806
    "AC" = -0.07017433399895806
    "AD" = 75.59999999999999
    "AF" = 75.59999999999999
    "AG" = -1841.5
    "AH" = 0.04396882633476795
    "AI" = 0.06972240154938669
    "AJ" = 9.314618108244813
    "AK" = 2.211215177713737
}

foreach ($row in 2, 3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
